# The "2024" sheet tracks a log of entries per month, with the most recent
# entry for each month-group kept at the top of its block. A new September
# entry was logged ("transfer share anyone axis" at 2024-09-05 16:03:14),
# which pushes all the existing entries in the "Others" group (and every
# row below it) down by one row.
#
# This is exactly what inserting a new row at the top of that block does:
# Excel shifts rows 30..66 down to 31..67 (growing the sheet's dimension
# from A1:Y66 to A1:Y67) and leaves a blank row 30 for the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new row above row 30 (the first data row of the "Others" group),
# pushing everything at/after row 30 down by one.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row with the new September log entry.
$ws.Range("R30").Value = "transfer share anyone axis"
$ws.Range("S30").Value = "2024-09-05 16:03:14"
